$d = $word.ActiveDocument

# 1. "ARAF and MAP2K1" -> "ARAF, MAP2K1"  (no field codes in this span, plain replace is safe)
$d.Content.Find.Execute("ARAF and MAP2K1", $true, $false, $false, $false, $false, $true, 1, $false, "ARAF, MAP2K1", 2) | Out-Null

# 2. ". " -> ", " right after the "1-9" citation field, before "NF1 and PTPN11"
#    (must not span the citation field itself, so find the literal text run just after it)
$r2 = $d.Content
$r2.Find.Execute(". NF1 and PTPN11") | Out-Null
$seg2 = $d.Range($r2.Start, $r2.Start + 2)
$seg2.Text = ", "

# 3. Remove the PTPN11 citation clause entirely, collapsing to a period
#    "NF1 and PTPN11 genes may also be mutated9 but are not targeted by this gene panel." -> "NF1 and PTPN11."
#    This span legitimately deletes the citation field (matches the diff).
$d.Content.Find.Execute("NF1 and PTPN11 genes may also be mutated9 but are not targeted by this gene panel.", $true, $false, $false, $false, $false, $true, 1, $false, "NF1 and PTPN11.", 2) | Out-Null

# 4. Drop the trailing clause after the NTRK1 citation field, replacing with a period.
#    Must not span the citation field, so only match the literal text run after it.
$r4 = $d.Content
$r4.Find.Execute(" however these are not detected by this assay.") | Out-Null
$r4.Text = "."

# 5. Remove stray _GoBack bookmark split in "cellular heterogeneity" (no text change, merges runs + drops bookmark)
$d.Content.Find.Execute("cellular heterogeneity", $true, $false, $false, $false, $false, $true, 1, $false, "cellular heterogeneity", 2) | Out-Null
